$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3979.2
$ws.Range("I76").Value = 3966
$ws.Range("K76").Value = 3966
$ws.Range("M76").Value = -3651

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3979.2
$ws.Range("I79").Value = 3966
$ws.Range("K79").Value = 3966
$ws.Range("M79").Value = -2874

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 18520880
$ws.Range("I106").Value = 23810848
$ws.Range("K106").Value = 23810848
$ws.Range("M106").Value = -23810217

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 599.2
$ws.Range("J107").Value = 499.5
$ws.Range("L107").Value = 499.5
$ws.Range("N107").Value = -4339.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 156560.7
$ws.Range("I132").Value = 336249.75
$ws.Range("K132").Value = 1008749.25
$ws.Range("M132").Value = -1006219.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6976.8857
$ws.Range("J138").Value = 8999.392
$ws.Range("L138").Value = 26998.176
$ws.Range("N138").Value = -37278.176

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 73694.44500000001
$ws.Range("J140").Value = 74416.664
$ws.Range("L140").Value = 74416.664
$ws.Range("N140").Value = -84776.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 720455.0600000001
$ws.Range("I2").Value = 1249145.9
$ws.Range("K2").Value = 1249145.9
$ws.Range("M2").Value = -1249032.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 285.7143
$ws.Range("J4").Value = 300
$ws.Range("L4").Value = 300
$ws.Range("N4").Value = -532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1995.6
$ws.Range("I45").Value = 2328
$ws.Range("J45").Value = 1497
$ws.Range("K45").Value = 2328
$ws.Range("L45").Value = 1497
$ws.Range("M45").Value = -1951
$ws.Range("N45").Value = -2251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 4000
$ws.Range("J50").Value = 3500
$ws.Range("L50").Value = 3500
$ws.Range("N50").Value = -4928

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 39270.848
$ws.Range("J88").Value = 72228.28999999999
$ws.Range("L88").Value = 72228.28999999999
$ws.Range("N88").Value = -73040.28999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 39270.848
$ws.Range("J91").Value = 72228.28999999999
$ws.Range("L91").Value = 72228.28999999999
$ws.Range("N91").Value = -75036.28999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1858155.1
$ws.Range("I110").Value = 2551838.2
$ws.Range("K110").Value = 2551838.2
$ws.Range("M110").Value = -2549793.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 720455.0600000001
$ws.Range("I116").Value = 1249145.9
$ws.Range("K116").Value = 1249145.9
$ws.Range("M116").Value = -1246851.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 110182.164
$ws.Range("J140").Value = 110182.164
$ws.Range("L140").Value = 110182.164
$ws.Range("N140").Value = -120542.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 720455.0600000001
$ws.Range("I3").Value = 1249145.9
$ws.Range("K3").Value = 1249145.9
$ws.Range("M3").Value = -1249031.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2650
$ws.Range("I105").Value = 3271.2
$ws.Range("J105").Value = 2411.077
$ws.Range("K105").Value = 3271.2
$ws.Range("L105").Value = 2411.077
$ws.Range("M105").Value = -1524.2
$ws.Range("N105").Value = -5905.077

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1688.762
$ws.Range("I107").Value = 2442.111
$ws.Range("J107").Value = 1123.75
$ws.Range("K107").Value = 2442.111
$ws.Range("L107").Value = 1123.75
$ws.Range("M107").Value = -522.1109999999999
$ws.Range("N107").Value = -4963.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 129858
$ws.Range("J140").Value = 129858
$ws.Range("L140").Value = 129858
$ws.Range("N140").Value = -140218

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1569.6428
$ws.Range("I16").Value = 1489.091
$ws.Range("J16").Value = 1865
$ws.Range("K16").Value = 1489.091
$ws.Range("L16").Value = 1865
$ws.Range("M16").Value = -1202.091
$ws.Range("N16").Value = -2439

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2273832.8
$ws.Range("I105").Value = 3247446.8
$ws.Range("K105").Value = 3247446.8
$ws.Range("M105").Value = -3245699.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 793173.4
$ws.Range("I107").Value = 1299731.6
$ws.Range("K107").Value = 1299731.6
$ws.Range("M107").Value = -1297811.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1569.6428
$ws.Range("I113").Value = 1489.091
$ws.Range("J113").Value = 1865
$ws.Range("K113").Value = 1489.091
$ws.Range("L113").Value = 1865
$ws.Range("M113").Value = 680.9090000000001
$ws.Range("N113").Value = -6205

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3249.3333
$ws.Range("I134").Value = 3030.5
$ws.Range("K134").Value = 9091.5
$ws.Range("M134").Value = -6556.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 91662.05
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 97180.11
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 97180.11
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -107540.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 258
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 216
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 648
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -3020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 258
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 216
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 1944
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -13800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 848886.4
$ws.Range("I122").Value = 1002820.25
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 3008460.75
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3006010.75
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 56727.25
$ws.Range("J123").Value = 56727.25
$ws.Range("L123").Value = 56727.25
$ws.Range("N123").Value = -61627.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9998.333000000001
$ws.Range("I132").Value = 9997.5
$ws.Range("K132").Value = 29992.5
$ws.Range("M132").Value = -27462.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3249.5625
$ws.Range("I61").Value = 3622.9092
$ws.Range("J61").Value = 2428.2
$ws.Range("K61").Value = 3622.9092
$ws.Range("L61").Value = 2428.2
$ws.Range("M61").Value = -3420.9092
$ws.Range("N61").Value = -2832.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3249.5625
$ws.Range("I113").Value = 3622.9092
$ws.Range("J113").Value = 2428.2
$ws.Range("K113").Value = 3622.9092
$ws.Range("L113").Value = 2428.2
$ws.Range("M113").Value = -1452.9092
$ws.Range("N113").Value = -6768.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1762.7567
$ws.Range("I107").Value = 1813.2174
$ws.Range("K107").Value = 5439.6522
$ws.Range("M107").Value = -3519.6522

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 606.43475
$ws.Range("I113").Value = 471.53333
$ws.Range("J113").Value = 859.375
$ws.Range("K113").Value = 1414.59999
$ws.Range("L113").Value = 2578.125
$ws.Range("M113").Value = 755.4000100000001
$ws.Range("N113").Value = -6918.125
